# Apply the "wind_direction_definitions" restructuring + header rename
# described in the commit "updating code to use .xlsx file.".
#
# 1. wind_direction_definitions (sheet1): swap columns so "direction" (text)
#    comes first and "angle" (number) comes second, drop the trailing
#    360 -> N wrap-around row (9 data rows instead of 10), and resize the
#    columns to fit the new, narrower content.
# 2. region_direction_parameters (sheet4): rename header row from
#    "Standard"/"Region" to "standard"/"wind_region".
# 3. Restore the view/selection state left behind by the edit: the
#    wind_direction_definitions tab ends up active/selected, and the
#    other two touched sheets keep a plain (non-active) selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. wind_direction_definitions
# ---------------------------------------------------------------------
$wsDir = $wb.Worksheets.Item("wind_direction_definitions")

# Drop the old 10th row (angle 360 -> N) - the new table only needs 9 rows.
$wsDir.Rows.Item(10).Delete()

$directions = @("N", "NE", "E", "SE", "S", "SW", "W", "NW")
$angles = @(0, 45, 90, 135, 180, 225, 270, 315)

$wsDir.Cells.Item(1, 1).Value = "direction"
$wsDir.Cells.Item(1, 2).Value = "angle"

for ($i = 0; $i -lt $directions.Count; $i++) {
    $r = $i + 2
    $wsDir.Cells.Item($r, 1).Value = $directions[$i]
    $wsDir.Cells.Item($r, 2).Value = $angles[$i]
}

# Narrower columns now that column A holds short direction codes and
# column B holds small angle numbers instead of the reverse.
$wsDir.Columns.Item(1).ColumnWidth = 8.166666666666666
$wsDir.Columns.Item(2).ColumnWidth = 5.022222222222222

# ---------------------------------------------------------------------
# 2. region_direction_parameters header rename
# ---------------------------------------------------------------------
$wsRegionDir = $wb.Worksheets.Item("region_direction_parameters")
$wsRegionDir.Cells.Item(1, 1).Value = "standard"
$wsRegionDir.Cells.Item(1, 2).Value = "wind_region"

# ---------------------------------------------------------------------
# 3. View / selection state
# ---------------------------------------------------------------------
$wsWindspeed = $wb.Worksheets.Item("region_windspeed_parameters")
$wsWindspeed.Activate()
$wsWindspeed.Range("I30").Select()

$wsRegionDir.Activate()
$wsRegionDir.Range("D2:D9").Select()

$wsDir.Activate()
$wsDir.Range("B9").Select()
